# The original workbook has a title/header row (row 1, a merged cell
# "question_template_per_category" spanning A1:B1) above the actual
# "Category"/"Question" table header. The commit removes that header
# line, so everything below shifts up by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:1").Delete()
